# This edit re-shuffles the per-row price/listing data (columns D, H, I, J, K, L, M, N, O, P)
# across rows 2-19 of the sheet, as if reassigning each market listing to a different
# date/row, per the weekly update. Column A, B, C, E, F, G, Q, R are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values of the affected columns for every data row (2-19)
# before any of them get overwritten, since values move between rows.
$orig = @{}
for ($r = 2; $r -le 19; $r++) {
    $rowData = @{}
    $rowData["D"] = $ws.Cells.Item($r, 4).Value2
    $rowData["H"] = $ws.Cells.Item($r, 8).Value2
    $rowData["I"] = $ws.Cells.Item($r, 9).Value2
    $rowData["J"] = $ws.Cells.Item($r, 10).Value2
    $rowData["K"] = $ws.Cells.Item($r, 11).Value2
    $rowData["L"] = $ws.Cells.Item($r, 12).Value2
    $rowData["M"] = $ws.Cells.Item($r, 13).Value2
    $rowData["N"] = $ws.Cells.Item($r, 14).Value2
    $rowData["O"] = $ws.Cells.Item($r, 15).Value2
    $rowData["P"] = $ws.Cells.Item($r, 16).Value2
    $orig[$r] = $rowData
}

# Mapping of target row -> source row (i.e. target row receives source row's data)
$mapping = @{
    2 = 19
    3 = 10
    4 = 13
    5 = 14
    6 = 11
    7 = 15
    8 = 17
    9 = 5
    10 = 2
    11 = 3
    12 = 9
    13 = 6
    14 = 4
    15 = 18
    16 = 8
    17 = 16
    18 = 7
    19 = 12
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $src = $orig[$sourceRow]
    $ws.Cells.Item($targetRow, 4).Value2 = $src["D"]
    $ws.Cells.Item($targetRow, 8).Value2 = $src["H"]
    $ws.Cells.Item($targetRow, 9).Value2 = $src["I"]
    $ws.Cells.Item($targetRow, 10).Value2 = $src["J"]
    $ws.Cells.Item($targetRow, 11).Value2 = $src["K"]
    $ws.Cells.Item($targetRow, 12).Value2 = $src["L"]
    $ws.Cells.Item($targetRow, 13).Value2 = $src["M"]
    $ws.Cells.Item($targetRow, 14).Value2 = $src["N"]
    $ws.Cells.Item($targetRow, 15).Value2 = $src["O"]
    $ws.Cells.Item($targetRow, 16).Value2 = $src["P"]
}
